$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new "school" column before the existing "detail" column ---
# (old B:detail, C:amount, D:date) -> (new B:school [blank for old rows],
#  C:detail, D:amount, E:date)
$ws.Columns("B:B").Insert()
$ws.Range("B2").Value = "school"

# --- Append the two new data rows (16 and 17) ---
$ws.Range("A16").Value = 12
$ws.Range("B16").Value = "โรงเรียนเทพสถิตวิทยา"
$ws.Range("C16").Value = "ค่าวัสดุการศึกษา"
$ws.Range("D16").Value = 2
# "2023-12-29" looks like a date, so Excel would normally auto-convert it to
# a date serial on direct assignment. Write it as a formula producing the
# text, then copy/paste-values so it lands as a literal string (matching
# the plain text values used in the rest of the "date" column).
$ws.Range("E16").Formula = "=""2023-12-29"""
$ws.Range("E16").Copy()
$ws.Range("E16").PasteSpecial(-4163)

$ws.Range("A17").Value = 13
$ws.Range("B17").Value = "โรงเรียนเขาดินพิทยารักษ์"
$ws.Range("C17").Value = "ค่าวัสดุการศึกษา"
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = "Created on 29-12-2023"

# --- Restore the selection Excel leaves active after these edits ---
$ws.Range("B3").Select() | Out-Null
